# EIA Table 2.10.A update: October 2016/2015 -> November 2016/2015 (2017-01-31 update)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: assign a literal text value to a cell without letting Excel's
# automatic type inference reinterpret a date-like string (e.g. "November
# 2016") as a date serial number. We briefly force a Text number format,
# write the value, then restore the exact original format code so the
# cell's style index is unchanged.
# ---------------------------------------------------------------------------
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.NumberFormat = $fmt
}

# ---------------------------------------------------------------------------
# Title / period labels
# ---------------------------------------------------------------------------
$ws.Range("A2").Value2 = "November 2016 and November 2015 (Thousand Tons)"

Set-TextValue "B5" "November 2016"
Set-TextValue "C5" "November 2015"
Set-TextValue "E5" "November 2016"
Set-TextValue "F5" "November 2015"
Set-TextValue "G5" "November 2016"
Set-TextValue "H5" "November 2015"
Set-TextValue "I5" "November 2016"
Set-TextValue "J5" "November 2015"
Set-TextValue "K5" "November 2016"
Set-TextValue "L5" "November 2015"

# ---------------------------------------------------------------------------
# Data updates (Middle Atlantic block - row 13)
# ---------------------------------------------------------------------------
$ws.Range("C13").Value2 = 4
$ws.Range("L13").Value2 = 4

# New Jersey - row 14 (values become "NM" - Not Meaningful)
$ws.Range("B14").Value2 = "NM"
$ws.Range("D14").Value2 = "NM"
$ws.Range("K14").Value2 = "NM"

# East North Central - row 17
$ws.Range("B17").Value2 = 59
$ws.Range("C17").Value2 = 115
$ws.Range("D17").Value2 = -0.49
$ws.Range("E17").Value2 = 24
$ws.Range("F17").Value2 = 64
$ws.Range("G17").Value2 = 31
$ws.Range("H17").Value2 = 47
$ws.Range("K17").Value2 = 4
$ws.Range("L17").Value2 = 4

# Indiana - row 19
$ws.Range("C19").Value2 = 55
$ws.Range("F19").Value2 = 55

# Michigan - row 20
$ws.Range("C20").Value2 = 9
$ws.Range("D20").Value2 = 1.65
$ws.Range("E20").Value2 = 20
$ws.Range("F20").Value2 = 5

# Ohio - row 21
$ws.Range("B21").Value2 = 31
$ws.Range("C21").Value2 = 45
$ws.Range("D21").Value2 = -0.31
$ws.Range("G21").Value2 = 31
$ws.Range("H21").Value2 = 45

# Wisconsin - row 22
$ws.Range("B22").Value2 = 4
$ws.Range("C22").Value2 = 6
$ws.Range("D22").Value2 = -0.22
$ws.Range("E22").Value2 = 4
$ws.Range("F22").Value2 = 4
$ws.Range("K22").Value2 = 0.36
$ws.Range("L22").Value2 = 2

# West North Central - row 23
$ws.Range("I23").Value2 = 0.08
$ws.Range("J23").Value2 = 0.21

# Iowa - row 24
$ws.Range("I24").Value2 = 0.08
$ws.Range("J24").Value2 = 0.21

# South Atlantic - row 31
$ws.Range("B31").Value2 = 18
$ws.Range("C31").Value2 = 35
$ws.Range("D31").Value2 = -0.49
$ws.Range("E31").Value2 = 16
$ws.Range("F31").Value2 = 32
$ws.Range("K31").Value2 = 2

# Florida - row 34
$ws.Range("B34").Value2 = 16
$ws.Range("C34").Value2 = 32
$ws.Range("D34").Value2 = -0.5
$ws.Range("E34").Value2 = 16
$ws.Range("F34").Value2 = 32

# Georgia - row 35
$ws.Range("B35").Value2 = 2
$ws.Range("D35").Value2 = -0.42
$ws.Range("K35").Value2 = 2

# East South Central - row 41
$ws.Range("B41").Value2 = 31
$ws.Range("C41").Value2 = 30
$ws.Range("D41").Value2 = 0.02
$ws.Range("E41").Value2 = 31
$ws.Range("F41").Value2 = 30

# Kentucky - row 43
$ws.Range("B43").Value2 = 31
$ws.Range("C43").Value2 = 30
$ws.Range("D43").Value2 = 0.02
$ws.Range("E43").Value2 = 31
$ws.Range("F43").Value2 = 30

# West South Central - row 46
$ws.Range("B46").Value2 = 175
$ws.Range("C46").Value2 = 59
$ws.Range("D46").Value2 = 1.96
$ws.Range("E46").Value2 = 169
$ws.Range("F46").Value2 = 52
$ws.Range("K46").Value2 = 6
$ws.Range("L46").Value2 = 7

# Louisiana - row 48
$ws.Range("B48").Value2 = 172
$ws.Range("C48").Value2 = 55
$ws.Range("D48").Value2 = 2.12
$ws.Range("E48").Value2 = 169
$ws.Range("F48").Value2 = 52
$ws.Range("L48").Value2 = 3

# Mountain - row 51
$ws.Range("B51").Value2 = 15
$ws.Range("C51").Value2 = 15
$ws.Range("D51").Value2 = 0.039
$ws.Range("G51").Value2 = 15
$ws.Range("H51").Value2 = 15

# Nevada - row 55
$ws.Range("B55").Value2 = 15
$ws.Range("C55").Value2 = 15
$ws.Range("D55").Value2 = 0.039
$ws.Range("G55").Value2 = 15
$ws.Range("H55").Value2 = 15

# U.S. Total - row 67
$ws.Range("B67").Value2 = 304
$ws.Range("C67").Value2 = 260
$ws.Range("D67").Value2 = 0.17
$ws.Range("E67").Value2 = 240
$ws.Range("F67").Value2 = 178
$ws.Range("G67").Value2 = 47
$ws.Range("H67").Value2 = 62
$ws.Range("I67").Value2 = 0.08
$ws.Range("J67").Value2 = 0.21
$ws.Range("K67").Value2 = 18
$ws.Range("L67").Value2 = 20
